# Generate Report for Handoff
#
# Re-running the handoff report writer updated the "Latest Handoff Datetime"
# for the 11a4a630-b121-4c2f-ba62-1ec45dea4bac.md entry on the "zh-cn"
# status sheet (row 5, column H) to reflect the new handoff timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("H5").Value = "2016-08-18 20:42:58"
